$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the touched cells keep their Text number format so values
# like "0.999" or "67.010.30" are not reinterpreted as numbers/dates.
$cells = @("D2","E2","D3","E3","E4","D5","E5","D6","E6","E7","D8","E8","D9","E9","D10","E10","E11","D12","E12","E13","D14","E14","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","E26","E27","E28","E29","E30","E31","E32","D33","E33","E34","D35","E35","E36","D37","E37","D38","D39","E39","D40","E40","E41","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","D48","E48","E50","E51")
foreach ($addr in $cells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = '67.010.30'
$ws.Range("E2").Value = '  +2.71%  '
$ws.Range("D3").Value = '3.093.62'
$ws.Range("E3").Value = '  +5.05%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '580.71'
$ws.Range("E5").Value = '  +2.01%  '
$ws.Range("D6").Value = '169.40'
$ws.Range("E6").Value = '  +6.47%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '3.087.96'
$ws.Range("E8").Value = '  +4.92%  '
$ws.Range("D9").Value = '0.524'
$ws.Range("E9").Value = '  +1.52%  '
$ws.Range("D10").Value = '6.62'
$ws.Range("E10").Value = '  -2.26%  '
$ws.Range("E11").Value = '  +3.85%  '
$ws.Range("D12").Value = '0.482'
$ws.Range("E12").Value = '  +4.83%  '
$ws.Range("E13").Value = '  +1.91%  '
$ws.Range("D14").Value = '36.42'
$ws.Range("E14").Value = '  +5.69%  '
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("D16").Value = '3.606.01'
$ws.Range("E16").Value = '  +5.07%  '
$ws.Range("D17").Value = '66.918.76'
$ws.Range("E17").Value = '  +2.51%  '
$ws.Range("D18").Value = '7.19'
$ws.Range("E18").Value = '  +2.41%  '
$ws.Range("D19").Value = '3.089.91'
$ws.Range("E19").Value = '  +5.00%  '
$ws.Range("D20").Value = '16.25'
$ws.Range("E20").Value = '  +4.80%  '
$ws.Range("D21").Value = '466.24'
$ws.Range("E21").Value = '  +4.79%  '
$ws.Range("D22").Value = '0.714'
$ws.Range("E22").Value = '  +2.87%  '
$ws.Range("D23").Value = '7.49'
$ws.Range("E23").Value = '  +2.70%  '
$ws.Range("D24").Value = '84.12'
$ws.Range("E24").Value = '  +1.99%  '
$ws.Range("D25").Value = '13.09'
$ws.Range("E25").Value = '  +7.59%  '
$ws.Range("E26").Value = '  +2.96%  '
$ws.Range("E27").Value = '  +0.35%  '
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("E29").Value = '  +0.21%  '
$ws.Range("E30").Value = '  +0.48%  '
$ws.Range("E31").Value = '  +3.54%  '
$ws.Range("E32").Value = '  +1.37%  '
$ws.Range("D33").Value = '28.30'
$ws.Range("E33").Value = '  +3.82%  '
$ws.Range("E34").Value = '  +2.29%  '
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("E36").Value = '  +3.47%  '
$ws.Range("D37").Value = '5.89'
$ws.Range("E37").Value = '  +2.41%  '
$ws.Range("D38").Value = '47.66'
$ws.Range("D39").Value = '2.09'
$ws.Range("E39").Value = '  +5.47%  '
$ws.Range("D40").Value = '50.60'
$ws.Range("E40").Value = '  +2.05%  '
$ws.Range("E41").Value = '  +4.54%  '
$ws.Range("E42").Value = '  +1.59%  '
$ws.Range("D43").Value = '8.66'
$ws.Range("E43").Value = '  +1.56%  '
$ws.Range("D44").Value = '2.79'
$ws.Range("E44").Value = '  -1.57%  '
$ws.Range("D45").Value = '0.0360'
$ws.Range("E45").Value = '  +1.99%  '
$ws.Range("D46").Value = '381.06'
$ws.Range("E46").Value = '  -0.33%  '
$ws.Range("D47").Value = '2.787.01'
$ws.Range("E47").Value = '  +3.41%  '
$ws.Range("D48").Value = '134.78'
$ws.Range("E48").Value = '  +0.64%  '
$ws.Range("E50").Value = '  +5.70%  '
$ws.Range("E51").Value = '  +1.14%  '
